$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data: an email address (hyperlinked) and the same
# "Short" label used in row 1.
$ws.Range("A2").Value = "dilan@utexas.edu"
$ws.Range("B2").Value = "Short"

# Turn the new email address into a mailto: hyperlink, mirroring A1.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:dilan@utexas.edu")

# Hyperlinks.Add() stamps its own ad-hoc style on the cell; restore the
# existing "Hyperlink" cell style (the one already used by A1) by copying
# A1's formatting onto A2 instead of leaving the newly minted one in place.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Move the active selection like the author's next step would have.
[void]$ws.Range("B3").Select()
